# Apply translation-sheet edits:
#  - Header row 3: English header label changes from "GB" to "ENG" and the
#    GB/SI columns (F/G) are swapped so Slovenian (SI) comes first.
#  - For every existing data row, the English (F) and Slovenian (G) values
#    were swapped (Slovenian now in F, English now in G).
#  - Two new rows are appended for "Relays active" / "Radius".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Scratch cell used to swap two cells while preserving their original cell
# type (so numeric-looking text like "0" stays text instead of becoming a
# real number) - Range.Copy keeps the source type/format, unlike assigning
# through .Value2 which re-infers the type from the string.
$scratch = $ws.Range("Z1")

function Swap-Cells($a, $b) {
    $a.Copy($scratch)
    $b.Copy($a)
    $scratch.Copy($b)
}

# --- Header row (row 3): F/G language headers ---
Swap-Cells $ws.Range("F3") $ws.Range("G3")
$ws.Range("G3").Value2 = "ENG"

# --- Swap English/Slovenian values (columns F and G) for existing rows ---
$swapRows = @(4, 10, 11, 14, 15, 16, 17, 18, 21, 22, 25, 26, 29, 35, 36, 37, 38, 39, 40, 41, 43, 52, 53, 54)

foreach ($r in $swapRows) {
    Swap-Cells $ws.Cells.Item($r, 6) $ws.Cells.Item($r, 7)
}

$scratch.ClearContents()

# --- New row 55: "Relays active" ---
$ws.Cells.Item(55, 2).Value2 = "SingleUseId57"
$ws.Cells.Item(55, 3).Value2 = "Default"
$ws.Cells.Item(55, 4).Value2 = "Left"
$ws.Cells.Item(55, 5).Value2 = "LTR"
$ws.Cells.Item(55, 6).Value2 = "Releji aktivni"
$ws.Cells.Item(55, 7).Value2 = "Relays active"

# --- New row 56: "Radius" ---
$ws.Cells.Item(56, 2).Value2 = "SingleUseId58"
$ws.Cells.Item(56, 3).Value2 = "Default"
$ws.Cells.Item(56, 4).Value2 = "Left"
$ws.Cells.Item(56, 5).Value2 = "LTR"
$ws.Cells.Item(56, 6).Value2 = "Radius"
$ws.Cells.Item(56, 7).Value2 = "Radius"
